# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets.
# Mirrors the upstream gh-pages data refresh commit (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# "展览" sheet (rows keyed by row number, column F)
$exhibitUpdates = @{
    2  = 1108
    3  = 422
    4  = 1515
    5  = 8795
    8  = 653
    9  = 296
    12 = 16
    13 = 3647
    16 = 84
    17 = 2577
    21 = 212
    22 = 2440
    23 = 71
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value = $exhibitUpdates[$row]
}

# "全部类型" sheet (rows keyed by row number, column F)
$allTypesUpdates = @{
    2  = 1108
    3  = 422
    4  = 1515
    5  = 8795
    8  = 653
    9  = 296
    12 = 16
    13 = 3647
    16 = 84
    17 = 2578
    21 = 212
    22 = 2440
    24 = 71
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
